$wb = $excel.ActiveWorkbook

# Sheet: ALC (index 1), Row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1710.9
$ws.Range("I40").Value = 1299.8
$ws.Range("J40").Value = 2122
$ws.Range("K40").Value = 1299.8
$ws.Range("L40").Value = 2122
$ws.Range("M40").Value = -1124.8
$ws.Range("N40").Value = -2472

# Sheet: ALC (index 1), Row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 9007.546
$ws.Range("I64").Value = 6000
$ws.Range("J64").Value = 9308.299999999999
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 9308.299999999999
$ws.Range("M64").Value = -5752
$ws.Range("N64").Value = -9804.299999999999

# Sheet: ALC (index 1), Row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 9007.546
$ws.Range("I67").Value = 6000
$ws.Range("J67").Value = 9308.299999999999
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 9308.299999999999
$ws.Range("M67").Value = -5142
$ws.Range("N67").Value = -11024.3

# Sheet: ALC (index 1), Row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 11743.866
$ws.Range("I113").Value = 6487.1113
$ws.Range("K113").Value = 6487.1113
$ws.Range("M113").Value = -3233.1113

# Sheet: ALC (index 1), Row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 4499.6665
$ws.Range("I125").Value = 4250
$ws.Range("J125").Value = 4999
$ws.Range("K125").Value = 38250
$ws.Range("L125").Value = 44991
$ws.Range("M125").Value = -35790
$ws.Range("N125").Value = -49911

# Sheet: ALC (index 1), Row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 99564.5
$ws.Range("J130").Value = 99564.5
$ws.Range("L130").Value = 99564.5
$ws.Range("N130").Value = -109604.5

# Sheet: ALC (index 1), Row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1832.8572
$ws.Range("I132").Value = 1385.0834
$ws.Range("K132").Value = 4155.2502
$ws.Range("M132").Value = -1625.2502

# Sheet: ALC (index 1), Row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1819.7
$ws.Range("I137").Value = 1949.8334
$ws.Range("K137").Value = 5849.5002
$ws.Range("M137").Value = -3299.5002

# Sheet: ARM (index 2), Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6861.232
$ws.Range("I32").Value = 3324.2341
$ws.Range("K32").Value = 3324.2341
$ws.Range("M32").Value = -3037.2341

# Sheet: ARM (index 2), Row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 50599.5
$ws.Range("I45").Value = 50599.5
$ws.Range("K45").Value = 50599.5
$ws.Range("M45").Value = -50222.5

# Sheet: ARM (index 2), Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3677.6619
$ws.Range("I61").Value = 2419.0513
$ws.Range("J61").Value = 5370.276
$ws.Range("K61").Value = 2419.0513
$ws.Range("L61").Value = 5370.276
$ws.Range("M61").Value = -2207.0513
$ws.Range("N61").Value = -5794.276

# Sheet: ARM (index 2), Row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2250.3333
$ws.Range("I88").Value = 1711.6
$ws.Range("K88").Value = 1711.6
$ws.Range("M88").Value = -1305.6

# Sheet: ARM (index 2), Row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2250.3333
$ws.Range("I91").Value = 1711.6
$ws.Range("K91").Value = 1711.6
$ws.Range("M91").Value = -307.5999999999999

# Sheet: ARM (index 2), Row 130
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 125095.22
$ws.Range("J130").Value = 125095.22
$ws.Range("L130").Value = 125095.22
$ws.Range("N130").Value = -135135.22

# Sheet: ARM (index 2), Row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 106282.25
$ws.Range("J133").Value = 106282.25
$ws.Range("L133").Value = 106282.25
$ws.Range("N133").Value = -111342.25

# Sheet: ARM (index 2), Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3677.6619
$ws.Range("I136").Value = 2419.0513
$ws.Range("J136").Value = 5370.276
$ws.Range("K136").Value = 7257.1539
$ws.Range("L136").Value = 16110.828
$ws.Range("M136").Value = -4707.1539
$ws.Range("N136").Value = -21210.828

# Sheet: BSM (index 3), Row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1517.75
$ws.Range("I105").Value = 1679
$ws.Range("J105").Value = 1249
$ws.Range("K105").Value = 1679
$ws.Range("L105").Value = 1249
$ws.Range("M105").Value = 68
$ws.Range("N105").Value = -4743

# Sheet: CRP (index 4), Row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9661.532999999999
$ws.Range("I62").Value = 9703
$ws.Range("K62").Value = 9703
$ws.Range("M62").Value = -9079

# Sheet: CRP (index 4), Row 64
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 15833.333
$ws.Range("J64").Value = 15833.333
$ws.Range("L64").Value = 15833.333
$ws.Range("N64").Value = -16329.333

# Sheet: CRP (index 4), Row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 9661.532999999999
$ws.Range("I65").Value = 9703
$ws.Range("K65").Value = 48515
$ws.Range("M65").Value = -45395

# Sheet: CRP (index 4), Row 67
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 15833.333
$ws.Range("J67").Value = 15833.333
$ws.Range("L67").Value = 15833.333
$ws.Range("N67").Value = -17549.333

# Sheet: CUL (index 5), Row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1050.5714
$ws.Range("I92").Value = 199.5
$ws.Range("J92").Value = 1192.4166
$ws.Range("K92").Value = 598.5
$ws.Range("L92").Value = 3577.2498
$ws.Range("M92").Value = 649.5
$ws.Range("N92").Value = -6073.2498

# Sheet: GSM (index 6), Row 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 20101956
$ws.Range("I14").Value = 25126782
$ws.Range("J14").Value = 2650
$ws.Range("K14").Value = 25126782
$ws.Range("L14").Value = 2650
$ws.Range("M14").Value = -25126614
$ws.Range("N14").Value = -2986

# Sheet: GSM (index 6), Row 22
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 32999.5
$ws.Range("J22").Value = 32999.5
$ws.Range("L22").Value = 32999.5
$ws.Range("N22").Value = -34057.5

# Sheet: GSM (index 6), Row 25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 7001.6665
$ws.Range("I25").Value = 6008
$ws.Range("J25").Value = 7498.5
$ws.Range("K25").Value = 6008
$ws.Range("L25").Value = 7498.5
$ws.Range("M25").Value = -5479
$ws.Range("N25").Value = -8556.5

# Sheet: GSM (index 6), Row 52
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 43979.6
$ws.Range("J52").Value = 44974.5
$ws.Range("L52").Value = 44974.5
$ws.Range("N52").Value = -45492.5

# Sheet: GSM (index 6), Row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3001.739
$ws.Range("I122").Value = 2457.1
$ws.Range("K122").Value = 7371.299999999999
$ws.Range("M122").Value = -4921.299999999999

# Sheet: LTW (index 7), Row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 55172.15
$ws.Range("I61").Value = 57707.527
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 57707.527
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -57505.527
$ws.Range("N61").Value = -7404

# Sheet: LTW (index 7), Row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 55172.15
$ws.Range("I113").Value = 57707.527
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 57707.527
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = -55537.527
$ws.Range("N113").Value = -11340

# Sheet: LTW (index 7), Row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3871.9524
$ws.Range("I132").Value = 2892.7058
$ws.Range("K132").Value = 8678.117400000001
$ws.Range("M132").Value = -6148.117400000001

# Sheet: WVR (index 8), Row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 511.22223
$ws.Range("I100").Value = 468.72223
$ws.Range("J100").Value = 596.2222
$ws.Range("K100").Value = 937.44446
$ws.Range("L100").Value = 1192.4444
$ws.Range("M100").Value = -396.44446
$ws.Range("N100").Value = -2274.4444
